# The "UserModule" results sheet already has a placeholder, pre-formatted
# (bold font + shaded fill) empty cell at E2 reserved for the test outcome.
# Record that the test passed by writing "PASSED" into it, keeping the
# bold / shaded formatting that the cell already carries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserModule")
$cell = $ws.Range("E2")

$cell.Value = "PASSED"

# Make sure the cell keeps/gets the bold Calibri 11 formatting used for the
# other status cells in this column (E3/E4 use the same look).
$cell.Font.Name = "Calibri"
$cell.Font.Size = 11
$cell.Font.Bold = $true
